$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'polyester athletic pants men',
    'compression for the knee',
    'volleyball youth',
    'knee sleeve with padding',
    'kneeling pad exercise',
    'baseball items',
    'boys cold leggings',
    'patella knee band',
    'leg sleeves for basketball youth boys',
    'mesh capri leggings',
    'mens black compression pants',
    'best knee pads construction',
    'paintball pants small',
    'below the knee',
    'football pads',
    'leg sleeves for basketball youth',
    'knee pads for working',
    'volleyball shorts longer length',
    'kneepads construction',
    'youth softball pants black',
    'mountain skin pants',
    'padded snowboarding shorts',
    'volleyball kneepads',
    'boys basketball gear',
    'baseball leg guard',
    'construction gel knee pads',
    'girls youth softball pants black',
    'running knee band',
    'recovery pants',
    'knee foam pad',
    'lacrosse shorts girls',
    'compression sleeve knee pads',
    'compression running tights men',
    'bjj shorts for men',
    'football shorts for men',
    'weight basketball',
    'professional work knee pads',
    'wrestling shorts',
    'mesh basketball shorts for men',
    'catchers gear leg guard',
    'gel knee pads construction',
    'acl knee',
    'knee sleeve wrestling',
    'spandex shorts men',
    'sheer test',
    'black leggings youth',
    'calf sleeve padded',
    'compression hip',
    'knee sleeve for wrestling',
    'compression knee for men',
    'work knee pads',
    'knees pads',
    'medium youth baseball pants',
    'knee compression sleeve with knee pad',
    'knee pads with gel',
    'football compression shorts youth',
    'pants mountain',
    'hex soccer',
    'knee guards mountain biking',
    'construction kneeling pad',
    'leg guards softball',
    'sort pants men',
    'knee pad work',
    'construction work knee pads',
    'knee sleeves football',
    'knee sleeve running men',
    'squat pants men',
    'down pants',
    'adult baseball pants black',
    'long basketball shorts',
    'knee padding',
    'shorts pad',
    'leg sleeves for men football',
    'black baseball pants youth',
    'baseball shorts for boys',
    'construction knee pads gel',
    'mountain biking knee pads',
    'mountain biking pads for men',
    'work kneepads',
    'construction knee pads for work',
    'knee pads for works',
    'compression knee sleeves for weightlifting',
    'knee pads for work',
    'compression shorts football',
    'under shorts for men',
    'soccer tights',
    'knee pad floor',
    'gel work knee pads',
    'knee pads for kneeling',
    'knees pads for work',
    'biking capris',
    'hip protector',
    'volleyball kneepads women',
    'compression shorts bjj',
    'basketball sleeve youth boys',
    'yoga tights',
    'padded sliding shorts women',
    'tight leggings',
    'catchers gear women',
    'softball gear for men',
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
